# Update cryptocurrency price/volume figures (scheduled data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "66.160.98"
$ws.Cells.Item(2, 5).Value = "  -0.10%  "
$ws.Cells.Item(3, 4).Value = "3.018.38"
$ws.Cells.Item(3, 5).Value = "  +0.05%  "
$ws.Cells.Item(4, 4).Value = "'1.00"
$ws.Cells.Item(4, 5).Value = "  +0.02%  "
$ws.Cells.Item(5, 4).Value = "'576.84"
$ws.Cells.Item(5, 5).Value = "  -1.47%  "
$ws.Cells.Item(6, 4).Value = "'168.37"
$ws.Cells.Item(6, 5).Value = "  +3.15%  "
$ws.Cells.Item(7, 4).Value = "'1.00"
$ws.Cells.Item(7, 5).Value = "  -0.02%  "
$ws.Cells.Item(8, 4).Value = "'0.521"
$ws.Cells.Item(8, 5).Value = "  +0.43%  "
$ws.Cells.Item(9, 4).Value = "3.013.53"
$ws.Cells.Item(9, 5).Value = "  +0.12%  "
$ws.Cells.Item(10, 4).Value = "'6.74"
$ws.Cells.Item(10, 5).Value = "  +0.05%  "
$ws.Cells.Item(11, 5).Value = "  -1.98%  "
$ws.Cells.Item(12, 4).Value = "'0.477"
$ws.Cells.Item(12, 5).Value = "  +4.24%  "
$ws.Cells.Item(13, 5).Value = "  -3.40%  "
$ws.Cells.Item(14, 4).Value = "'37.06"
$ws.Cells.Item(14, 5).Value = "  +6.51%  "
$ws.Cells.Item(15, 5).Value = "  -0.44%  "
$ws.Cells.Item(16, 4).Value = "66.151.26"
$ws.Cells.Item(16, 5).Value = "  -0.02%  "
$ws.Cells.Item(17, 4).Value = "3.521.64"
$ws.Cells.Item(17, 5).Value = "  +0.07%  "
$ws.Cells.Item(18, 4).Value = "'7.27"
$ws.Cells.Item(18, 5).Value = "  +4.58%  "
$ws.Cells.Item(19, 4).Value = "3.025.62"
$ws.Cells.Item(19, 5).Value = "  +0.25%  "
$ws.Cells.Item(20, 4).Value = "'16.20"
$ws.Cells.Item(20, 5).Value = "  +16.19%  "
$ws.Cells.Item(21, 4).Value = "'468.85"
$ws.Cells.Item(21, 5).Value = "  +2.20%  "
$ws.Cells.Item(22, 4).Value = "'0.707"
$ws.Cells.Item(22, 5).Value = "  +2.45%  "
$ws.Cells.Item(23, 4).Value = "'7.48"
$ws.Cells.Item(23, 5).Value = "  +0.98%  "
$ws.Cells.Item(24, 4).Value = "'83.51"
$ws.Cells.Item(24, 5).Value = "  +1.30%  "
$ws.Cells.Item(25, 4).Value = "'12.82"
$ws.Cells.Item(25, 5).Value = "  +3.38%  "
$ws.Cells.Item(26, 4).Value = "'2.30"
$ws.Cells.Item(26, 5).Value = "  -0.59%  "
$ws.Cells.Item(27, 4).Value = "'10.19"
$ws.Cells.Item(27, 5).Value = "  -3.67%  "
$ws.Cells.Item(28, 5).Value = "  +0.21%  "
$ws.Cells.Item(29, 4).Value = "'8.48"
$ws.Cells.Item(29, 5).Value = "  +4.83%  "
$ws.Cells.Item(30, 4).Value = "'2.46"
$ws.Cells.Item(30, 5).Value = "  +3.61%  "
$ws.Cells.Item(31, 4).Value = "'2.62"
$ws.Cells.Item(31, 5).Value = "  +0.08%  "
$ws.Cells.Item(32, 5).Value = "  -4.00%  "
$ws.Cells.Item(33, 5).Value = "  +6.07%  "
$ws.Cells.Item(34, 4).Value = "'28.13"
$ws.Cells.Item(34, 5).Value = "  +2.80%  "
$ws.Cells.Item(35, 4).Value = "'1.00"
$ws.Cells.Item(35, 5).Value = "  +0.11%  "
$ws.Cells.Item(36, 2).Value = "Filecoin"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Cells.Item(36, 4).Value = "'5.85"
$ws.Cells.Item(36, 5).Value = "  +0.29%  "
$ws.Cells.Item(37, 2).Value = "Mantle"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Cells.Item(37, 4).Value = "'0.989"
$ws.Cells.Item(37, 5).Value = "  -0.50%  "
$ws.Cells.Item(38, 4).Value = "'48.24"
$ws.Cells.Item(38, 5).Value = "  +10.49%  "
$ws.Cells.Item(39, 4).Value = "'2.05"
$ws.Cells.Item(39, 5).Value = "  -7.02%  "
$ws.Cells.Item(40, 2).Value = "OKB"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Cells.Item(40, 4).Value = "'49.51"
$ws.Cells.Item(40, 5).Value = "  -1.11%  "
$ws.Cells.Item(41, 2).Value = "TheGraph"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Cells.Item(41, 4).Value = "'0.313"
$ws.Cells.Item(41, 5).Value = "  +1.12%  "
$ws.Cells.Item(42, 5).Value = "  -4.61%  "
$ws.Cells.Item(43, 5).Value = "  -1.08%  "
$ws.Cells.Item(44, 4).Value = "'8.64"
$ws.Cells.Item(44, 5).Value = "  +2.00%  "
$ws.Cells.Item(45, 5).Value = "  -0.95%  "
$ws.Cells.Item(46, 4).Value = "'382.71"
$ws.Cells.Item(46, 5).Value = "  -3.55%  "
$ws.Cells.Item(47, 4).Value = "2.728.81"
$ws.Cells.Item(47, 5).Value = "  -2.71%  "
$ws.Cells.Item(48, 4).Value = "'133.71"
$ws.Cells.Item(48, 5).Value = "  -0.12%  "
$ws.Cells.Item(49, 5).Value = "  -0.01%  "
$ws.Cells.Item(50, 4).Value = "'24.79"
$ws.Cells.Item(50, 5).Value = "  +3.17%  "
$ws.Cells.Item(51, 5).Value = "  +3.66%  "
